# Applies the "next step adding figs to Supplementary materials" edit:
#  1. Bump the fixed "Date" footer placeholder text (slide master + all
#     slide layouts) from 20-Dec-20 to 21-Dec-20.
#  2. On slide 4 (the Fig1/Fig2 timeline slide), widen/relabel the
#     "Expansion" callout -> "Main Analysis" (bold) and widen/bold the
#     "Phase III" callout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: look through a shape container's Shapes collection for a date
# placeholder whose current text is the old fixed date, and update it.
# ---------------------------------------------------------------------
function Update-FixedDateShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "20-Dec-20") {
                $sh.TextFrame.TextRange.Text = "21-Dec-20"
            }
        }
    }
}

# ---------------------------------------------------------------------
# Helper: find a shape by its exact Name within a slide's Shapes.
# ---------------------------------------------------------------------
function Find-ShapeByName($container, $name) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# 1. Slide master date placeholder.
$master = $p.SlideMaster
Update-FixedDateShape $master

# 1b. Every custom (slide) layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-FixedDateShape $master.CustomLayouts.Item($li)
}

# 2. Slide 4 timeline callouts.
$slide4 = $p.Slides.Item(4)

# "Expansion" -> "Main Analysis": move/widen box, bold, relabel.
$expansionBox = Find-ShapeByName $slide4 "TextBox 42"
if ($expansionBox -ne $null) {
    $expansionBox.Left = 181.41944181889764
    $expansionBox.Width = 118.5012628425197
    $expansionBox.TextFrame.TextRange.Text = "Main Analysis"
    $expansionBox.TextFrame.TextRange.Font.Bold = 1
}

# "Phase III": move/widen box, bold (keep existing italic + text).
$phaseIIIBox = Find-ShapeByName $slide4 "TextBox 118"
if ($phaseIIIBox -ne $null) {
    $phaseIIIBox.Left = 60.2348051496063
    $phaseIIIBox.Width = 101.36086714173229
    $phaseIIIBox.TextFrame.TextRange.Font.Bold = 1
}

Write-Host "edit.ps1 applied"
